# BillHubTestdata.xlsx — finish the "TaxCode" sheet and make it the active tab.
#
# The TaxCode sheet gets a small header table:
#   Tax Code | Tax Percent | Description
#   ZC        | 12          | Valid test data
#   INVALIDTC | invalidpercentage | Invalid test data
#   (blank)   | 13          | Update test data
#
# and becomes the workbook's active/selected sheet (previously "BA" was
# the selected tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaxCode")

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Tax Code"
$ws.Range("B1").Value = "Tax Percent"
$ws.Range("C1").Value = "Description"

# --- Data rows (write order reproduces the author's original entry
# order / shared-string table layout) ----------------------------------
$ws.Range("A3").Value = "INVALIDTC"
$ws.Range("B3").Value = "invalidpercentage"
$ws.Range("C2").Value = "Valid test data"
$ws.Range("C3").Value = "Invalid test data"
$ws.Range("C4").Value = "Update test data"
$ws.Range("A2").Value = "ZC"
$ws.Range("B2").Value = 12
$ws.Range("B4").Value = 13

# --- Header formatting --------------------------------------------------
$ws.Range("A1:C1").Font.Bold = $true

# --- Column widths (approximate the author's manual resize) -----------
$ws.Columns.Item(1).ColumnWidth = 14.6
$ws.Columns.Item(2).ColumnWidth = 15.1
$ws.Columns.Item(3).ColumnWidth = 13.3

# --- Print setup --------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Make TaxCode the active sheet / selection -------------------------
$ws.Activate()
$ws.Range("A2").Select()
